$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column Q, shifting nothing (it's the first empty column to the
# right of the existing table) - this clones the per-row cell styles from
# column P into column Q for every row that has data, exactly like Excel
# does when you insert a column.
$ws.Columns("Q:Q").Insert(-4161)   # xlShiftToRight

# The insert also stamped a (style-only) Q15 cell because row 15 previously
# had a lone P15 cell; the source workbook does NOT have a Q15 cell, so
# remove it again (xlShiftUp - nothing below it, so this cleanly drops the
# cell without touching the rest of column Q).
$ws.Cells.Item(15, 17).Delete(-4162)   # xlShiftUp

# Fill in the new 2020 data column (Q) with its header + per-region values.
$ws.Cells.Item(4, 17).Value = 2020

$ws.Cells.Item(5, 17).Value = 38.6
$ws.Cells.Item(6, 17).Value = 42.4
$ws.Cells.Item(7, 17).Value = 53.2
$ws.Cells.Item(8, 17).Value = 90.6
$ws.Cells.Item(9, 17).Value = 52.6
$ws.Cells.Item(10, 17).Value = 24.5
$ws.Cells.Item(11, 17).Value = 69.1
$ws.Cells.Item(12, 17).Value = 32.2
$ws.Cells.Item(13, 17).Value = 19.1
$ws.Cells.Item(14, 17).Value = 25.2

# Match the author's final selection state.
$ws.Range("R27").Select()

$wb.Save()
